$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.638.53"
$ws.Range("E2").Value = "  +1.59%  "

# Row 3
$ws.Range("D3").Value = "1.903.00"
$ws.Range("E3").Value = "  +2.61%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.45%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.19"
$ws.Range("E5").Value = "  +1.47%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.48%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4796"
$ws.Range("E7").Value = "  +0.52%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2834"
$ws.Range("E8").Value = "  +0.71%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06544"
$ws.Range("E9").Value = "  +0.62%  "

# Row 10
$ws.Range("D10").Value = "1.958.45"
$ws.Range("E10").Value = "  +4.77%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07493"
$ws.Range("E11").Value = "  +1.97%  "

# Row 12
$ws.Range("E12").Value = "  +1.94%  "

# Row 13
$ws.Range("E13").Value = "  -0.75%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.10"
$ws.Range("E14").Value = "  +1.11%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6688"
$ws.Range("E15").Value = "  +3.69%  "

# Row 16
$ws.Range("D16").Value = "30.617.05"
$ws.Range("E16").Value = "  +1.63%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.31"
$ws.Range("E17").Value = "  +0.77%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9993"
$ws.Range("E18").Value = "  -0.56%  "

# Row 19
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.204.93"
$ws.Range("E19").Value = "  +3.69%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007600"
$ws.Range("E20").Value = "  +0.02%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "227.57"
$ws.Range("E21").Value = "  +5.46%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.324"
$ws.Range("E22").Value = "  +1.11%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.51%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.215"
$ws.Range("E24").Value = "  +1.86%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.06"
$ws.Range("E25").Value = "  +2.55%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.282"
$ws.Range("E26").Value = "  -0.18%  "

# Row 27
$ws.Range("E27").Value = "  +0.35%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.960"
$ws.Range("E28").Value = "  +2.99%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.405"
$ws.Range("E29").Value = "  -1.57%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09834"
$ws.Range("E30").Value = "  +7.67%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.351"
$ws.Range("E31").Value = "  +2.84%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.026"
$ws.Range("E32").Value = "  +1.44%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05062"
$ws.Range("E33").Value = "  +0.82%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.230"
$ws.Range("E34").Value = "  +8.79%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7538"
$ws.Range("E35").Value = "  +1.66%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.712"
$ws.Range("E36").Value = "  +0.48%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01877"
$ws.Range("E37").Value = "  +2.87%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.659"
$ws.Range("E38").Value = "  +2.01%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9184"
$ws.Range("E39").Value = "  +1.14%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.080"
$ws.Range("E40").Value = "  +1.71%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "107.37"
$ws.Range("E41").Value = "  +0.30%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.845"
$ws.Range("E42").Value = "  -1.09%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4297"
$ws.Range("E43").Value = "  +1.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.007"
$ws.Range("E44").Value = "  +0.35%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.416"
$ws.Range("E45").Value = "  +0.44%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.53"
$ws.Range("E46").Value = "  +0.63%  "

# Row 47
$ws.Range("E47").Value = "  -2.67%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.051"
$ws.Range("E48").Value = "  +2.40%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.483"
$ws.Range("E49").Value = "  -4.79%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.83"
$ws.Range("E50").Value = "  -1.22%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05670"
$ws.Range("E51").Value = "  -0.73%  "
